# Update "想去人数" (attendee interest count) figures on the "展览" and
# "全部类型" worksheets to reflect freshly generated data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1816
    $ws.Range("F3").Value = 8214
    $ws.Range("F5").Value = 314
}
